$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = -0.2407239916940589
$ws.Range("C2").Value = 0.009033579491229321
$ws.Range("D2").Value = -0.2423913045098955
$ws.Range("E2").Value = -0.2390552569459463

$ws.Range("B3").Value = -0.5803439289158613
$ws.Range("C3").Value = 0.03274360758043506
$ws.Range("D3").Value = -0.5845858438557359
$ws.Range("E3").Value = -0.5760702868730736

$ws.Range("B4").Value = -0.1859378220230946
$ws.Range("C4").Value = 0.005112071127402038
$ws.Range("D4").Value = -0.1869049750396768
$ws.Range("E4").Value = -0.1849703085686199

$ws.Range("B5").Value = -0.2407239916940589
$ws.Range("C5").Value = 0.009033579491229321
$ws.Range("D5").Value = -0.2423913045098955
$ws.Range("E5").Value = -0.2390552569459463

$ws.Range("B6").Value = -0.4513711966981849
$ws.Range("C6").Value = 0.02166622749986542
$ws.Range("D6").Value = -0.4547466340866235
$ws.Range("E6").Value = -0.4479827925341042

$ws.Range("B7").Value = -0.1549396964772941
$ws.Range("C7").Value = 0.003963964878086172
$ws.Range("D7").Value = -0.1556978947484505
$ws.Range("E7").Value = -0.1541813156418506
